# Swap the Approved/Rejected verdicts (and ReasonToReject) between
# row 23 and row 24 of the Test-Cases sheet, and update the active
# cell selection to reflect the cell the user was editing (I23:J23).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 23 was "Approved" with no rejection reason; it becomes "Rejected"
# with reason "Nil".
$ws.Range("I23").Value = "Rejected"
$ws.Range("J23").Value = "Nil"

# Row 24 was "Rejected" with reason "Nil"; it becomes "Approved" with
# no rejection reason (cell cleared).
$ws.Range("I24").Value = "Approved"
$ws.Range("J24").Value = ""

# Reflect the new selection/active cell shown in the sheet view.
$ws.Range("I23:J23").Select() | Out-Null
